$d = $word.ActiveDocument

# The edit inserts the clause " MEDIANTE LA" right after
# "...ACREDITEN LA POSESION DEL PREDIO" (and before the pre-existing
# trailing space), turning:
#   "...ACREDITEN LA POSESIÓN DEL PREDIO FIRMA ENTRE LAS PARTES."
# into:
#   "...ACREDITEN LA POSESIÓN DEL PREDIO MEDIANTE LA FIRMA ENTRE LAS PARTES."
#
# The phrase "ACREDITEN LA POSESIÓN DEL PREDIO" (followed by a space) also
# occurs a few other times in the document with different continuations, so
# we disambiguate by searching for a longer string that also covers the
# start of the following run ("FIRMA ENTRE LAS PARTES"), which only occurs
# once in the whole document.
$anchorText = "ACREDITEN LA POSESIÓN DEL PREDIO"
$uniqueAnchor = $anchorText + " FIRMA ENTRE LAS PARTES"

$searchRange = $d.Content
$found = $searchRange.Find.Execute($uniqueAnchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Position right after "...ACREDITEN LA POSESIÓN DEL PREDIO" (i.e.
    # right before the trailing space that precedes "FIRMA...").
    $insertAt = $searchRange.Start + $anchorText.Length

    $insertionPoint = $d.Range($insertAt, $insertAt)
    $insertionPoint.InsertAfter(" MEDIANTE LA")
} else {
    Write-Host "WARNING: target anchor text not found; document left unchanged."
}
